$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Motor pin labels (rows 12-15) now filled in with new motor names
$ws.Range("E12").Value = "Motor L2"
$ws.Range("E13").Value = "Motor R1"
$ws.Range("E14").Value = "Motor L1"
$ws.Range("E15").Value = "Motor R2"

# Rows 23-24 cleared (duplicate motor labels removed)
$ws.Range("E23").Value = ""
$ws.Range("E24").Value = ""

# Rows 29-30: ultrasonic front sensor labels + I/O direction
$ws.Range("E29").Value = "F TRIG"
$ws.Range("F29").Value = "OUT"
$ws.Range("E30").Value = "F ECHO"
$ws.Range("F30").Value = "IN"

# Rows 33-36: add OUT direction labels, rename row 36 to Back HL
$ws.Range("F33").Value = "OUT"
$ws.Range("F34").Value = "OUT"
$ws.Range("F35").Value = "OUT"
$ws.Range("E36").Value = "Back HL"
$ws.Range("F36").Value = "OUT"

# Rows 37-38 cleared (duplicate motor labels removed)
$ws.Range("E37").Value = ""
$ws.Range("E38").Value = ""

# Rows 39-40: ultrasonic back sensor labels + I/O direction
$ws.Range("E39").Value = "B TRIG"
$ws.Range("F39").Value = "OUT"
$ws.Range("E40").Value = "B ECHO"
$ws.Range("F40").Value = "IN"

# Update sheet view (scroll position / selection)
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J32").Select()
